$wb = $excel.ActiveWorkbook

# Add a new worksheet at the end of the workbook and name it "nr_studies"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "nr_studies"

# Header row
$ws.Range("A1").Value = "outcome"
$ws.Range("B1").Value = "moderator_design"
$ws.Range("C1").Value = "n_effect_sizes"
$ws.Range("D1").Value = "k_studies"

$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

# Data rows
$data = @(
    @("NS", "Experimental (non-randomized)", 15, 2),
    @("NS", "Cross-sectional", 651, 65),
    @("NS", "Longitudinal", 38, 11),
    @("NS", "Cross-lagged", 3, 1),
    @("NS", "Experimental (RCT)", 13, 1),
    @("NT", "Cross-sectional", 356, 42),
    @("NT", "Longitudinal", 14, 6),
    @("NT", "Experimental (non-randomized)", 9, 2),
    @("NT", "Cross-lagged", 2, 1),
    @("NS", "Experimental (non-randomized", 2, 1)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

Write-Host "Added sheet nr_studies with $($data.Count) data rows"
